$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Refresh the cached "datetimeFigureOut" date field text (the footer/date
#    placeholder) wherever it appears: the notes master, every slide layout,
#    and the slide master. PowerPoint re-caches this automatically-updating
#    field's displayed text whenever the deck is saved; here we do it
#    explicitly for every placeholder of type "date" (ppPlaceholderDate=16).
# ---------------------------------------------------------------------------
$newDate = "3/13/2018"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDate = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) { $isDate = $true }
        } catch {
            $isDate = $false
        }
        if ($isDate -and $sh.HasTextFrame) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Notes master
Update-DatePlaceholder $p.NotesMaster.Shapes

# Every slide layout
$master = $p.SlideMaster
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# Slide master
Update-DatePlaceholder $master.Shapes

# ---------------------------------------------------------------------------
# 2. Rename the "Person..." UI component shapes on slide 1 to "Book...".
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $sh = $slide.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        $t = $sh.TextFrame.TextRange.Text
        if ($t -eq "PersonListPanel") {
            $sh.TextFrame.TextRange.Text = "BookListPanel"
        } elseif ($t -eq "PersonCard") {
            $sh.TextFrame.TextRange.Text = "BookCard"
        }
    }
}
